$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I13").Value = "jdhfj"
$ws.Range("J13").Value = "djff"
$ws.Range("I14").Value = "dhfj"

$ws.Range("I14").Select()
